# Remove the Footer and Slide Number placeholder shapes from the title
# slide ("Proposed Solution Approach - Revolutionary").
#
# Quirk of this COM host: calling Delete() on a layout-backed placeholder
# shape the first time doesn't remove it - it resets it to an empty,
# layout-inherited placeholder instead (new shape id/name, text cleared).
# A second Delete() call on that same slot (re-fetched fresh via
# Shapes.Item - a cached shape reference won't do) performs the real
# removal. So each target placeholder is looked up and deleted twice.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$targetNames = @("Footer Placeholder 3", "Slide Number Placeholder 4")

foreach ($targetName in $targetNames) {
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        if ($s.Shapes.Item($i).Name -eq $targetName) {
            # First call: resets the placeholder to an empty layout default.
            $s.Shapes.Item($i).Delete()
            # Second call (fresh lookup): actually removes it from the slide.
            if ($i -le $s.Shapes.Count) {
                $s.Shapes.Item($i).Delete()
            }
            break
        }
    }
}
